# Added week 1 hours: Liam Connors logged 8 hours in Week 1 (row 4, col B).
# The Total column (C) auto-recalculates via its existing cumulative-sum
# formulas, so we only need to set the new input cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8

# Move/restore the active selection to B5, matching the saved cursor
# position recorded in the worksheet view.
$ws.Range("B5").Select()
